# Refresh the crypto market snapshot: latest Price (D) and Volume(1h) (E) pulls,
# plus two rows (33/34) whose ranking swapped places (EthereumClassic <-> NEARProtocol).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, even when it parses as a number
# (e.g. "1.00", "65.00"), without leaving the cells number format/style altered -
# the source data models price/volume columns as plain text.
function Set-TextValue($range, $text) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $savedStyle
}

$ws.Range("D2").Value = "64.059.01"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "2.762.64"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue $ws.Range("D5") "576.59"
$ws.Range("E5").Value = "  -1.53%  "
Set-TextValue $ws.Range("D6") "159.33"
$ws.Range("E6").Value = "  -1.33%  "
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  +0.19%  "
Set-TextValue $ws.Range("D8") "0.601"
$ws.Range("E8").Value = "  -3.32%  "
$ws.Range("E9").Value = "  -3.61%  "
$ws.Range("E10").Value = "  +3.43%  "
Set-TextValue $ws.Range("D11") "5.80"
$ws.Range("E11").Value = "  -14.73%  "
Set-TextValue $ws.Range("D12") "0.385"
$ws.Range("E12").Value = "  -3.28%  "
$ws.Range("D13").Value = "3.250.14"
$ws.Range("E13").Value = "  -0.42%  "
Set-TextValue $ws.Range("D14") "26.91"
$ws.Range("D15").Value = "63.683.47"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("E16").Value = "  -5.09%  "
$ws.Range("D17").Value = "2.766.95"
$ws.Range("E17").Value = "  -0.68%  "
Set-TextValue $ws.Range("D18") "12.10"
$ws.Range("E18").Value = "  -1.97%  "
$ws.Range("E19").Value = "  -3.32%  "
Set-TextValue $ws.Range("D20") "356.10"
$ws.Range("E20").Value = "  -3.20%  "
$ws.Range("E21").Value = "  -5.56%  "
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("E23").Value = "  -6.33%  "
Set-TextValue $ws.Range("D24") "65.00"
$ws.Range("E24").Value = "  -3.52%  "
$ws.Range("E25").Value = "  -3.96%  "
$ws.Range("E26").Value = "  -2.21%  "
Set-TextValue $ws.Range("D27") "1.00"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").Value = "0.0₃0906"
$ws.Range("E28").Value = "  -6.49%  "
Set-TextValue $ws.Range("D29") "7.30"
$ws.Range("E29").Value = "  -0.89%  "
Set-TextValue $ws.Range("D30") "1.95"
$ws.Range("E30").Value = "  -3.91%  "
Set-TextValue $ws.Range("D31") "1.26"
$ws.Range("E31").Value = "  -0.02%  "
Set-TextValue $ws.Range("D32") "170.35"
$ws.Range("E32").Value = "  -1.38%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D33") "20.16"
$ws.Range("E33").Value = "  -3.33%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D34") "4.91"
$ws.Range("E34").Value = "  -3.46%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("E37").Value = "  -1.93%  "
$ws.Range("E38").Value = "  -3.32%  "
Set-TextValue $ws.Range("D39") "349.57"
$ws.Range("E39").Value = "  +1.99%  "
Set-TextValue $ws.Range("D40") "6.26"
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("E41").Value = "  -2.87%  "
Set-TextValue $ws.Range("D42") "39.08"
$ws.Range("E42").Value = "  -2.26%  "
Set-TextValue $ws.Range("D43") "21.46"
$ws.Range("E43").Value = "  -4.78%  "
Set-TextValue $ws.Range("D44") "21.76"
Set-TextValue $ws.Range("D45") "0.0586"
$ws.Range("E45").Value = "  -3.91%  "
Set-TextValue $ws.Range("D46") "138.07"
$ws.Range("E46").Value = "  -0.67%  "
Set-TextValue $ws.Range("D47") "0.631"
$ws.Range("E47").Value = "  -3.46%  "
$ws.Range("E48").Value = "  -3.03%  "
$ws.Range("E50").Value = "  +0.17%  "
Set-TextValue $ws.Range("D51") "11.04"
$ws.Range("E51").Value = "  -0.02%  "
